$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D4 was stored as the shared string "484"; Excel re-entered it as a
# true number, so the shared-strings table loses that entry on save. ---
$ws.Range("D4").Value = 484

# --- New year columns E:H (2020-2023) added alongside the existing D
# (2019) column; copy the existing column D number formats/styles across
# for each row so the new cells pick up the same formatting. ---
$ws.Range("D3:D6").Copy() | Out-Null
$ws.Range("E3:H3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Header year values
$ws.Range("E3").Value = 2020
$ws.Range("F3").Value = 2021
$ws.Range("G3").Value = 2022
$ws.Range("H3").Value = 2023

# Row 4: "Number of local governments" stayed flat at 484 for every year
$ws.Range("E4").Value = 484
$ws.Range("F4").Value = 484
$ws.Range("G4").Value = 484
$ws.Range("H4").Value = 484

# Row 5: proportion (%) series
$ws.Range("E5").Value = 13.2
$ws.Range("F5").Value = 21.5
$ws.Range("G5").Value = 34.5
$ws.Range("H5").Value = 40.53

# Row 6: count series
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 104
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 169

# --- Cosmetic tweaks that came along with the edit: narrower A:C columns,
# taller header row (to fit the now-longer wrapped header text), and the
# stale F6 cell-selection marker cleared back to the top-left cell. ---
$ws.Columns.Item(1).ColumnWidth = 39.25
$ws.Columns.Item(2).ColumnWidth = 39.25
$ws.Columns.Item(3).ColumnWidth = 39.25

$ws.Rows.Item(1).RowHeight = 79.5

$ws.Range("A1").Select() | Out-Null
